$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("A1").Value = "from"
$ws.Range("A2").Value = "from"
